# evidence submission for tasks A1-A6 (dankuzone w/ DAIC)

$wb = $excel.ActiveWorkbook

# --- Info sheet: just move the selection, no data changes ---
$wsInfo = $wb.Worksheets.Item("Info")
$wsInfo.Range("G7").Select()

# --- A1: TxHash / ClassID ---
$wsA1 = $wb.Worksheets.Item("A1")
$wsA1.Range("A2").Value = "6779BEE088F3AFD8AE23BF07FA7ED89B061021A310F055CD6DEB164168A14221"
$wsA1.Range("B2").Value = "saikoronotori"
$wsA1.Range("A2").Select()

# --- A2: TxHash / ClassID / NFTID (two rows of evidence) ---
$wsA2 = $wb.Worksheets.Item("A2")
$wsA2.Range("A2").Value = "71BE265098FCB57E7DFE83559E4FC350B742E7E4BA915FFEAF799295FAE0223E"
$wsA2.Range("B2").Value = "saikoronotori"
$wsA2.Range("C2").Value = "saikoronotori01"
$wsA2.Range("A3").Value = "DA82E9D40F97180067FB64D7F9BF84441B7F1DC222F91DE3E76D427651B8FC07"
$wsA2.Range("B3").Value = "saikoronotori"
$wsA2.Range("C3").Value = "saikoronotori02"
$wsA2.Range("A3").Select()

# --- A3: TxHash / ClassID / NFTID / ChainID ---
$wsA3 = $wb.Worksheets.Item("A3")
$wsA3.Range("A2").Value = "BCC192B1D57C88A0598B001CF389C4542504196758E5E6C8CC6ACFF3EE6DD73F"
$wsA3.Range("D2").Value = "uni-6"
$wsA3.Range("B2").Value = "juno14acs4qq74005wyucgv06lqzje8er88753rn47fdvqpf4gwxf92ysj9hz8n"
$wsA3.Range("C2").Value = "saikoronotori01"
$wsA3.Range("B10").Select()

# --- A4: TxHash / ClassID / NFTID / ChainID ---
$wsA4 = $wb.Worksheets.Item("A4")
$wsA4.Range("A2").Value = "998B1DF93777728F3B76B5FBA0213AA25B76BEA4296946A2BD57F76A9085E4A5"
$wsA4.Range("D2").Value = "`tgon-flixnet-1"
$wsA4.Range("B2").Value = "ibc/CCDF0E008EB8940349BAD859F198A7A07A52BEC85CC5FBC53AC80F81FDEBBB85"
$wsA4.Range("C2").Value = "saikoronotori02"
$wsA4.Range("A2").Select()

# --- A5: TxHash / ClassID / NFTID / ChainID ---
$wsA5 = $wb.Worksheets.Item("A5")
$wsA5.Range("A2").Value = "3538EFA85341F38E0EC645FF5D7396AF33A4CC3FB0CA92E50D5605BCFC781F69"
$wsA5.Range("D2").Value = "uni-6"
$wsA5.Range("B2").Value = "juno14acs4qq74005wyucgv06lqzje8er88753rn47fdvqpf4gwxf92ysj9hz8n"
$wsA5.Range("C2").Value = "saikoronotori01"
$wsA5.Range("A2").Select()

# --- A6: TxHash / ClassID / NFTID / ChainID (left as the final active sheet) ---
$wsA6 = $wb.Worksheets.Item("A6")
$wsA6.Range("D2").Value = "gon-flixnet-1"
$wsA6.Range("A2").Value = "657FACAED665C6C58BD896D8EFE307491980209D029400BF08A1A0B012B067C4"
$wsA6.Range("B2").Value = "ibc/CCDF0E008EB8940349BAD859F198A7A07A52BEC85CC5FBC53AC80F81FDEBBB85"
$wsA6.Range("C2").Value = "saikoronotori02"
$wsA6.Range("C9").Select()
